# Add 9 new match-result rows (342-350, date 2025-07-19) to the "Partidos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

# Each entry: fecha, jugador, equipo, posicion, goles, autogoles, arquero, goles_recibidos, tarjetas_amarillas, tarjetas_rojas, asistencias, Penales_Atajados
$rows = @(
    @("7/19/2025", "Edwin Casas",            "Amarillo", "Arquero",        0, 0, $true,  0, 0, 0, 0, 0),
    @("7/19/2025", "Alexander Uribe",        "Amarillo", "Mediocampista",  3, 0, $false, 0, 0, 0, 2, 0),
    @("7/19/2025", "David Fernando Velasco", "Amarillo", "Delantero",      1, 0, $false, 0, 0, 0, 0, 0),
    @("7/19/2025", "Juan David Espinal",     "Amarillo", "Mediocampista",  1, 0, $false, 0, 0, 0, 0, 0),
    @("7/19/2025", "Julio Cesar Castaño",    "Amarillo", "Mediocampista",  1, 0, $false, 0, 0, 0, 0, 0),
    @("7/19/2025", "Juan Carlos Otero",      "Amarillo", "Mediocampista",  1, 0, $false, 0, 0, 0, 0, 0),
    @("7/19/2025", "Arnul David Narvaez",    "Amarillo", "Delantero",      0, 0, $false, 0, 0, 0, 2, 0),
    @("7/19/2025", "Fabian Caicedo",         "Azul",     "Arquero",        0, 0, $true,  7, 0, 0, 0, 0),
    @("7/19/2025", "Oscar Eduardo Herrera",  "Azul",     "Mediocampista",  0, 0, $false, 0, 0, 1, 0, 0)
)

$startRow = 342
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    $ws.Cells.Item($r, 10).Value = $data[9]
    $ws.Cells.Item($r, 11).Value = $data[10]
    $ws.Cells.Item($r, 12).Value = $data[11]
}

# Selection after the edit lands on the row right after the new data.
$ws.Range("A352").Select() | Out-Null
